$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New "Section General-Purpose I/O" table (rows 12-16), mirroring the
# existing "Section Interrupt" table's formatting by copying formats
# from the analogous rows above (row1 = header style, rows2-5 = the
# alternating data-row styles) and then overwriting values.
# ---------------------------------------------------------------------

# Row 12: header row - same look as row 1 (fill2/border, centered)
$ws.Range("A1:C1").Copy()
$ws.Range("A12:C12").PasteSpecial(-4122)
$ws.Range("A12:C12").HorizontalAlignment = -4108
$ws.Range("A12").Value = "Sense Encode[6]"
$ws.Range("B12").Value = "Sense Encode[5]"
$ws.Range("C12").Value = "Sense Description"

# Row 13: same look as row 2
$ws.Range("A2:C2").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)
$ws.Range("A13").Value = 0
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = "Rising Edge"

# Row 14: same look as row 3
$ws.Range("A3:C3").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)
$ws.Range("A14").Value = 0
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "Falling Edge"

# Row 15: same look as row 4
$ws.Range("A4:C4").Copy()
$ws.Range("A15:C15").PasteSpecial(-4122)
$ws.Range("A15").Value = 1
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = "Change Level"

# Row 16: same look as row 5
$ws.Range("A5:C5").Copy()
$ws.Range("A16:C16").PasteSpecial(-4122)
$ws.Range("A16").Value = 1
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = "Reserved"

# ---------------------------------------------------------------------
# Leftover column-width metadata on F:H (columns 6-8), as in the diff.
# ---------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 18
$ws.Columns.Item(7).ColumnWidth = 18.5
$ws.Columns.Item(8).ColumnWidth = 20.833333333333332

# ---------------------------------------------------------------------
# Selection moves to B10; print orientation explicitly set to portrait.
# ---------------------------------------------------------------------
$ws.Range("B10").Select()
$ws.PageSetup.Orientation = 1
